$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'aave', 'aave', 4862456468, 15211264.30116868),
    @(3, 'ada', 'cardano', 30655930678, 36362934086.75214),
    @(4, 'algo', 'algorand', 2346196853, 8706896055.327736),
    @(5, 'ape', 'apecoin', 508773481, 799455492),
    @(6, 'apt', 'aptos', 3260757572, 684332972.1450868),
    @(7, 'arb', 'arbitrum', 2441731436, 5150239630),
    @(8, 'asr', 'as-roma-fan-token', 48687191, 7745096),
    @(9, 'avax', 'avalanche-2', 10376800230, 422276596.0335201),
    @(10, 'axs', 'axie-infinity', 420270732, 166320222.4290032),
    @(11, 'bch', 'bitcoin-cash', 12009670414, 19910696.77165078),
    @(12, 'bgb', 'bitget-token', 5055144911, 1139992035.978791),
    @(13, 'bonk', 'bonk', 2033551954, 77419592329436.58),
    @(14, 'bnt', 'bancor', 94767919, 113836542.5460683),
    @(15, 'cro', 'crypto-com-chain', 5282395400, 32384192415.06367),
    @(16, 'crv', 'curve-dao-token', 1385270155, 1385810267),
    @(17, 'cvx', 'convex-finance', 372419702, 81990703.72425415),
    @(18, 'doge', 'dogecoin', 35452563636, 150419716383.7052),
    @(19, 'dot', 'polkadot', 6317026150, 1522267060),
    @(20, 'dydx', 'dydx-chain', 516350093, 773348704.4859886),
    @(21, 'ens', 'ethereum-name-service', 968898216, 33165585.05450796),
    @(22, 'etc', 'ethereum-classic', 3619378642, 153091726.9519364),
    @(23, 'fet', 'fetch-ai', 1885862249, 2604959126.672),
    @(24, 'fil', 'filecoin', 1794763554, 683864884),
    @(25, 'gala', 'gala', 817148692, 45613440228.21294),
    @(26, 'hbar', 'hedera-hashgraph', 10984361437, 42392676005.9379),
    @(27, 'inj', 'injective-protocol', 1448962586, 97727220.33),
    @(28, 'jasmy', 'jasmycoin', 880853907, 48419999999.3058),
    @(29, 'kas', 'kaspa', 2525455590, 26363807936.34678),
    @(30, 'kava', 'kava', 425472961, 1082853140),
    @(31, 'ldo', 'lido-dao', 1391369525, 895888862.5921584),
    @(32, 'leo', 'leo-2', 1508205, 999890689.477294),
    @(33, 'link', 'chainlink', 16087346903, 678099970.4527868),
    @(34, 'mana', 'decentraland', 590322260, 1919188956.588888),
    @(35, 'near', 'near', 3499959779, 1246897382),
    @(36, 'ondo', 'ondo-finance', 3320027594, 3159107529),
    @(37, 'paxg', 'pax-gold', 949628420, 284182.589),
    @(38, 'pendle', 'pendle', 912572593, 167445463.8527003),
    @(39, 'pepe', 'pepe', 5112432040, 420690000000000),
    @(40, 'sand', 'the-sandbox', 734572750, 2445857126.223322),
    @(41, 'sei', 'sei-network', 1877534811, 5781805555),
    @(42, 'shib', 'shiba-inu', 7977124285, 589246036829115.8),
    @(43, 'sui', 'sui', 13551437650, 3511924479.569998),
    @(44, 'uni', 'uniswap', 6901381154, 600483073.71),
    @(45, 'xlm', 'stellar', 13877106318, 31298889990.9876),
    @(46, 'xmr', 'monero', 4556243690, 18446744.07370955),
    @(47, 'xrp', 'ripple', 190750152823, 59308385925)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
